# Refresh the crypto price/volume table (columns D = Price, E = Volume(1h))
# for rows 2-51 with the latest scraped values. Cells whose new text would
# otherwise be auto-parsed as a plain number (e.g. "0.9987", "306.60") are
# forced to stay text: set NumberFormat "@" before the write, then reset the
# cell back to the "Normal" style afterwards so no stray style survives and
# the cell matches the original inlineStr/plain-text representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.214.70"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "1.872.36"
$ws.Range("E3").Value = "  -1.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5174"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3746"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07161"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8940"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("D12").Value = "1.872.70"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07528"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9987"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008538"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9987"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "27.255.78"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("E21").Value = "  -2.37%  "
$ws.Range("D22").Value = "2.105.53"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.479"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.832"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.45%  "
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.093"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.685"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.689"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09245"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05136"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.083"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.162"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7289"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.142"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02034"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.533"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.076"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5318"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.539"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.349"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1477"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4642"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9983"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("E48").Value = "  -4.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.566"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.25%  "
